# Landscaping Data.xlsx - append new observation rows (352-358) collected
# for 2025-06-29 (serial 45837), matching the "Artifical Cut off" data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ row = 352; A = 45837; B = "Flowering";     C = "Large";  D = 71; E = 86; G = 0.28; H = 0.3;  I = "Yes"; J = 2; K = "Neutral"; L = 8; M = 0.61; N = 71; O = 30.01; P = 6; Q = 0.1; R = 9.9; S = 48; T = 0 },
    @{ row = 353; A = 45837; B = "Nonflowering";  C = "Medium"; D = 71; E = 86; G = 0.28; H = 0.3;  I = "Yes"; J = 3; K = "Bright";  L = 8; M = 0.61; N = 71; O = 30.01; P = 6; Q = 0.1; R = 9.9; S = 48; T = 0 },
    @{ row = 354; A = 45837; B = "Nonflowering";  C = "Small";  D = 71; E = 86; G = 0.28; H = 0.3;  I = "Yes"; J = 3; K = "Bright";  L = 8; M = 0.61; N = 71; O = 30.01; P = 6; Q = 0.1; R = 9.9; S = 48; T = 0 },
    @{ row = 355; A = 45837; B = "Nonflowering";  C = "Medium"; D = 71; E = 86; G = 0.28; H = 0.6;  I = "Yes"; J = 3; K = "Neutral"; L = 8; M = 0.61; N = 71; O = 30.01; P = 6; Q = 0.1; R = 9.9; S = 48; T = 0 },
    @{ row = 356; A = 45837; B = "Nonflowering";  C = "Medium"; D = 71; E = 86; G = 0.28; H = 0.5;  I = "Yes"; J = 3; K = "Neutral"; L = 8; M = 0.61; N = 71; O = 30.01; P = 6; Q = 0.1; R = 9.9; S = 48; T = 0 },
    @{ row = 357; A = 45837; B = "Nonflowering";  C = "Large";  D = 71; E = 86; G = 0.28; H = 0.7;  I = "Yes"; J = 4; K = "Dark";    L = 8; M = 0.61; N = 71; O = 30.01; P = 6; Q = 0.1; R = 9.9; S = 48; T = 0 },
    @{ row = 358; A = 45837; B = "Tree";          C = "Medium"; D = 71; E = 86; G = 0.28; H = 1.25; I = "Yes"; J = 1; K = "Neutral"; L = 8; M = 0.61; N = 71; O = 30.01; P = 6; Q = 0.1; R = 9.9; S = 48; T = 0 }
)

foreach ($r in $newRows) {
    $n = $r.row
    $ws.Range("A$n").Value2 = $r.A
    $ws.Range("A$n").NumberFormat = $ws.Range("A351").NumberFormat
    $ws.Range("B$n").Value2 = $r.B
    $ws.Range("C$n").Value2 = $r.C
    $ws.Range("D$n").Value2 = $r.D
    $ws.Range("E$n").Value2 = $r.E
    # column F is the ABS(D-E) temp-diff formula, filled in as a shared
    # formula across the whole new block after the data is in place
    $ws.Range("G$n").Value2 = $r.G
    $ws.Range("H$n").Value2 = $r.H
    $ws.Range("I$n").Value2 = $r.I
    $ws.Range("J$n").Value2 = $r.J
    $ws.Range("K$n").Value2 = $r.K
    $ws.Range("L$n").Value2 = $r.L
    $ws.Range("M$n").Value2 = $r.M
    $ws.Range("N$n").Value2 = $r.N
    $ws.Range("O$n").Value2 = $r.O
    $ws.Range("P$n").Value2 = $r.P
    $ws.Range("Q$n").Value2 = $r.Q
    $ws.Range("R$n").Value2 = $r.R
    $ws.Range("S$n").Value2 = $r.S
    $ws.Range("T$n").Value2 = $r.T
}

# Fill-down the Temp_Diff formula (same pattern as the existing F347:F351
# shared formula) across the newly added rows.
$ws.Range("F352:F358").Formula = "=ABS(D352-E352)"

# The used range now reaches row 358; scroll the view down to the new rows
# and land the selection where the author left off.
$excel.ActiveWindow.ScrollRow = 334
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P359").Select() | Out-Null
